$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.134.30"
$ws.Range("E2").Value = "  +5.61%  "

$ws.Range("D3").Value = "1.921.63"
$ws.Range("E3").Value = "  +2.50%  "

$ws.Range("E4").Value = "  -0.66%  "

$ws.Range("D5").Value = "'331.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.97%  "

$ws.Range("E6").Value = "  -0.65%  "

$ws.Range("D7").Value = "'0.5218"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "'0.4092"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.77%  "

$ws.Range("D9").Value = "'0.08539"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.07%  "

$ws.Range("D10").Value = "'43.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.80%  "

$ws.Range("E11").Value = "  +2.15%  "

$ws.Range("D12").Value = "'22.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.04%  "

$ws.Range("D13").Value = "'6.424"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.17%  "

$ws.Range("D14").Value = "1.924.02"
$ws.Range("E14").Value = "  +2.53%  "

$ws.Range("D15").Value = "'7.449"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.26%  "

$ws.Range("E16").Value = "  -0.61%  "

$ws.Range("D17").Value = "'96.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.30%  "

$ws.Range("E18").Value = "  +0.99%  "

$ws.Range("D19").Value = "'0.06712"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.27%  "

$ws.Range("D20").Value = "'18.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.34%  "

$ws.Range("E21").Value = "  -0.65%  "

$ws.Range("E22").Value = "  +2.23%  "

$ws.Range("D23").Value = "30.142.54"
$ws.Range("E23").Value = "  +5.55%  "

$ws.Range("D24").Value = "'11.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.85%  "

$ws.Range("D25").Value = "'2.223"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").Value = "2.144.82"
$ws.Range("E26").Value = "  +2.36%  "

$ws.Range("D27").Value = "'21.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.66%  "

$ws.Range("D28").Value = "'160.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.79%  "

$ws.Range("E29").Value = "  +1.85%  "

$ws.Range("D30").Value = "'129.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.30%  "

$ws.Range("D31").Value = "'1.083"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.99%  "

$ws.Range("E32").Value = "  +1.52%  "

$ws.Range("D33").Value = "'6.114"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.45%  "

$ws.Range("D34").Value = "'3.647"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.91%  "

$ws.Range("D35").Value = "'0.02512"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.47%  "

$ws.Range("D36").Value = "'0.06615"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.89%  "

$ws.Range("D37").Value = "'0.2219"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.38%  "

$ws.Range("D38").Value = "'5.247"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.33%  "

$ws.Range("D39").Value = "'1.241"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.92%  "

$ws.Range("D40").Value = "'8.963"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.36%  "

$ws.Range("D41").Value = "'0.6552"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.64%  "

$ws.Range("D42").Value = "'11.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.68%  "

$ws.Range("D43").Value = "'1.248"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.11%  "

$ws.Range("D44").Value = "'0.6187"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.96%  "

$ws.Range("E45").Value = "  +2.24%  "

$ws.Range("D46").Value = "'3.780"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.45%  "

$ws.Range("D47").Value = "'2.098"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.73%  "

$ws.Range("E48").Value = "  +2.69%  "

$ws.Range("D49").Value = "'125.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.52%  "

$ws.Range("D50").Value = "'1.165"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.98%  "

$ws.Range("D51").Value = "'79.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.52%  "
